# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos.xlsx sheet, matching the GitHub Actions data refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.317.16'
$ws.Range("E2").Value = '  +1.43%  '

$ws.Range("D3").Value = '1.833.25'
$ws.Range("E3").Value = '  +0.89%  '

$ws.Range("E4").Value = '  +0.90%  '

$ws.Range("D5").Value = '''315.08'
$ws.Range("E5").Value = '  +1.90%  '

$ws.Range("D7").Value = '''0.4739'
$ws.Range("E7").Value = '  +1.83%  '

$ws.Range("D8").Value = '''0.3689'
$ws.Range("E8").Value = '  +0.95%  '

$ws.Range("D9").Value = '''0.07443'
$ws.Range("E9").Value = '  +1.04%  '

$ws.Range("D10").Value = '''0.8858'
$ws.Range("E10").Value = '  +1.78%  '

$ws.Range("D11").Value = '''20.49'
$ws.Range("E11").Value = '  +1.16%  '

$ws.Range("D12").Value = '1.871.84'
$ws.Range("E12").Value = '  +3.73%  '

$ws.Range("D13").Value = '''0.07340'
$ws.Range("E13").Value = '  +3.21%  '

$ws.Range("D14").Value = '''5.438'
$ws.Range("E14").Value = '  +0.95%  '

$ws.Range("D15").Value = '''94.05'
$ws.Range("E15").Value = '  +3.07%  '

$ws.Range("D16").Value = '''6.570'
$ws.Range("E16").Value = '  +0.82%  '

$ws.Range("E17").Value = '  +0.62%  '

$ws.Range("D18").Value = '''0.000008798'
$ws.Range("E18").Value = '  +1.23%  '

$ws.Range("D20").Value = '27.532.41'
$ws.Range("E20").Value = '  +2.14%  '

$ws.Range("E21").Value = '  +1.03%  '

$ws.Range("D22").Value = '''5.288'
$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("E23").Value = '  +0.78%  '

$ws.Range("D24").Value = '2.096.32'
$ws.Range("E24").Value = '  +2.42%  '

$ws.Range("D25").Value = '''1.895'
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").Value = '''152.08'
$ws.Range("E26").Value = '  +0.67%  '

$ws.Range("D27").Value = '''18.68'
$ws.Range("E27").Value = '  +1.49%  '

$ws.Range("D28").Value = '''2.146'
$ws.Range("E28").Value = '  +0.37%  '

$ws.Range("D29").Value = '''5.232'
$ws.Range("E29").Value = '  -0.43%  '

$ws.Range("D30").Value = '''117.17'
$ws.Range("E30").Value = '  +0.72%  '

$ws.Range("E31").Value = '  +1.23%  '

$ws.Range("D32").Value = '''0.7501'
$ws.Range("E32").Value = '  -1.19%  '

$ws.Range("D33").Value = '''1.175'
$ws.Range("E33").Value = '  +0.84%  '

$ws.Range("D34").Value = '''4.547'
$ws.Range("E34").Value = '  +1.51%  '

$ws.Range("D35").Value = '''2.947'
$ws.Range("E35").Value = '  +1.75%  '

$ws.Range("E36").Value = '  +0.95%  '

$ws.Range("D37").Value = '''1.095'
$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("D38").Value = '''0.05345'
$ws.Range("E38").Value = '  +1.03%  '

$ws.Range("D39").Value = '''0.01954'
$ws.Range("E39").Value = '  +0.34%  '

$ws.Range("D40").Value = '''2.974'
$ws.Range("E40").Value = '  -0.16%  '

$ws.Range("D41").Value = '''2.404'
$ws.Range("E41").Value = '  +3.20%  '

$ws.Range("D42").Value = '''7.232'
$ws.Range("E42").Value = '  +1.06%  '

$ws.Range("D43").Value = '''0.5310'
$ws.Range("E43").Value = '  +0.35%  '

$ws.Range("D44").Value = '''0.1660'
$ws.Range("E44").Value = '  +0.27%  '

$ws.Range("D45").Value = '''8.483'
$ws.Range("E45").Value = '  +0.60%  '

$ws.Range("D46").Value = '''0.4940'
$ws.Range("E46").Value = '  +1.83%  '

$ws.Range("D47").Value = '''10.57'
$ws.Range("E47").Value = '  +1.56%  '

$ws.Range("E48").Value = '  +0.91%  '

$ws.Range("D49").Value = '''105.01'
$ws.Range("E49").Value = '  +1.70%  '

$ws.Range("D50").Value = '''1.672'
$ws.Range("E50").Value = '  +0.67%  '

$ws.Range("D51").Value = '''0.06298'
$ws.Range("E51").Value = '  +0.06%  '

